$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 205
$ws.Range("I4").Value = 205
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 205
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -91
$ws.Range("N4").ClearContents()

$ws.Range("H18").Value = 2833.3333
$ws.Range("I18").Value = 2833.3333
$ws.Range("K18").Value = 2833.3333
$ws.Range("M18").Value = -2549.3333

$ws.Range("H21").Value = 3687.375
$ws.Range("I21").Value = 3999.8572
$ws.Range("K21").Value = 3999.8572
$ws.Range("M21").Value = -3531.8572

$ws.Range("H23").Value = 3687.375
$ws.Range("I23").Value = 3999.8572
$ws.Range("K23").Value = 3999.8572
$ws.Range("M23").Value = -3765.8572

$ws.Range("H86").Value = 3050
$ws.Range("I86").Value = 3060
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3060
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1937
$ws.Range("N86").Value = -5246

$ws.Range("H88").Value = 1958.7
$ws.Range("J88").Value = 1401.5555
$ws.Range("L88").Value = 1401.5555
$ws.Range("N88").Value = -2213.5555

$ws.Range("H89").Value = 3050
$ws.Range("I89").Value = 3060
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 15300
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -9684
$ws.Range("N89").Value = -26232

$ws.Range("H91").Value = 1958.7
$ws.Range("J91").Value = 1401.5555
$ws.Range("L91").Value = 1401.5555
$ws.Range("N91").Value = -4209.5555

$ws.Range("H106").Value = 11799.04
$ws.Range("I106").Value = 3548
$ws.Range("J106").Value = 14404.632
$ws.Range("K106").Value = 3548
$ws.Range("L106").Value = 14404.632
$ws.Range("M106").Value = -2917
$ws.Range("N106").Value = -15666.632

$ws.Range("H107").Value = 4381.24
$ws.Range("I107").Value = 4230.727
$ws.Range("K107").Value = 4230.727
$ws.Range("M107").Value = -2310.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 744.2222
$ws.Range("I2").Value = 737.375
$ws.Range("K2").Value = 737.375
$ws.Range("M2").Value = -624.375

$ws.Range("H32").Value = 112169.36
$ws.Range("I32").Value = 119111.53
$ws.Range("K32").Value = 119111.53
$ws.Range("M32").Value = -118824.53

$ws.Range("H61").Value = 6440.3
$ws.Range("I61").Value = 1962.9286
$ws.Range("J61").Value = 16887.5
$ws.Range("K61").Value = 1962.9286
$ws.Range("L61").Value = 16887.5
$ws.Range("M61").Value = -1750.9286
$ws.Range("N61").Value = -17311.5

$ws.Range("H116").Value = 744.2222
$ws.Range("I116").Value = 737.375
$ws.Range("K116").Value = 737.375
$ws.Range("M116").Value = 1556.625

$ws.Range("H136").Value = 6440.3
$ws.Range("I136").Value = 1962.9286
$ws.Range("J136").Value = 16887.5
$ws.Range("K136").Value = 5888.7858
$ws.Range("L136").Value = 50662.5
$ws.Range("M136").Value = -3338.7858
$ws.Range("N136").Value = -55762.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 744.2222
$ws.Range("I3").Value = 737.375
$ws.Range("K3").Value = 737.375
$ws.Range("M3").Value = -623.375

$ws.Range("H20").Value = 49833.383
$ws.Range("I20").Value = 64836.188
$ws.Range("K20").Value = 64836.188
$ws.Range("M20").Value = -64589.188

$ws.Range("H86").Value = 3619.9333
$ws.Range("I86").Value = 3254.5454
$ws.Range("J86").Value = 4624.75
$ws.Range("K86").Value = 3254.5454
$ws.Range("L86").Value = 4624.75
$ws.Range("M86").Value = -2131.5454
$ws.Range("N86").Value = -6870.75

$ws.Range("H89").Value = 3619.9333
$ws.Range("I89").Value = 3254.5454
$ws.Range("J89").Value = 4624.75
$ws.Range("K89").Value = 16272.727
$ws.Range("L89").Value = 23123.75
$ws.Range("M89").Value = -10656.727
$ws.Range("N89").Value = -34355.75

$ws.Range("H134").Value = 11002.385
$ws.Range("I134").Value = 3303.2
$ws.Range("K134").Value = 9909.599999999999
$ws.Range("M134").Value = -7374.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 62500
$ws.Range("I3").Value = 60000
$ws.Range("J3").Value = 80000
$ws.Range("K3").Value = 60000
$ws.Range("L3").Value = 80000
$ws.Range("M3").Value = -59887
$ws.Range("N3").Value = -80226

$ws.Range("H31").Value = 2310.1333
$ws.Range("I31").Value = 2172.1428
$ws.Range("J31").Value = 2537.4119
$ws.Range("K31").Value = 2172.1428
$ws.Range("L31").Value = 2537.4119
$ws.Range("M31").Value = -1877.1428
$ws.Range("N31").Value = -3127.4119

$ws.Range("H34").Value = 2310.1333
$ws.Range("I34").Value = 2172.1428
$ws.Range("J34").Value = 2537.4119
$ws.Range("K34").Value = 2172.1428
$ws.Range("L34").Value = 2537.4119
$ws.Range("M34").Value = -1970.1428
$ws.Range("N34").Value = -2941.4119

$ws.Range("H94").Value = 6176.905
$ws.Range("J94").Value = 1402.5
$ws.Range("L94").Value = 1402.5
$ws.Range("N94").Value = -2304.5

$ws.Range("H122").Value = 6463.763
$ws.Range("I122").Value = 1609.5927
$ws.Range("K122").Value = 4828.7781
$ws.Range("M122").Value = -2378.7781

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5444
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30224

$ws.Range("H55").Value = 49416108
$ws.Range("J55").Value = 5644.615
$ws.Range("L55").Value = 16933.845
$ws.Range("N55").Value = -17287.845

$ws.Range("H109").Value = 2669.6667
$ws.Range("I109").Value = 1033.7142
$ws.Range("J109").Value = 4960
$ws.Range("K109").Value = 3101.1426
$ws.Range("L109").Value = 14880
$ws.Range("M109").Value = -2061.1426
$ws.Range("N109").Value = -16960

$ws.Range("H123").Value = 11150.7
$ws.Range("I123").Value = 8202
$ws.Range("J123").Value = 14099.4
$ws.Range("K123").Value = 24606
$ws.Range("L123").Value = 42298.2
$ws.Range("M123").Value = -22156
$ws.Range("N123").Value = -47198.2

$ws.Range("H131").Value = 3046.4583
$ws.Range("J131").Value = 3371.238
$ws.Range("L131").Value = 10113.714
$ws.Range("N131").Value = -20193.714

$ws.Range("H132").Value = 813.93335
$ws.Range("I132").Value = 655.5714
$ws.Range("K132").Value = 5900.1426
$ws.Range("M132").Value = -3370.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 925
$ws.Range("J4").Value = 925
$ws.Range("L4").Value = 925
$ws.Range("N4").Value = -1149

$ws.Range("H122").Value = 62463.53
$ws.Range("I122").Value = 103413.2
$ws.Range("J122").Value = 3964
$ws.Range("K122").Value = 310239.6
$ws.Range("L122").Value = 11892
$ws.Range("M122").Value = -307789.6
$ws.Range("N122").Value = -16792

$ws.Range("H132").Value = 9078.473
$ws.Range("I132").Value = 11218.167
$ws.Range("J132").Value = 4799.0835
$ws.Range("K132").Value = 33654.501
$ws.Range("L132").Value = 14397.2505
$ws.Range("M132").Value = -31124.501
$ws.Range("N132").Value = -19457.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5001
$ws.Range("J4").Value = 5001
$ws.Range("L4").Value = 5001
$ws.Range("N4").Value = -5227

$ws.Range("H22").Value = 1796.4333
$ws.Range("I22").Value = 599.5
$ws.Range("K22").Value = 599.5
$ws.Range("M22").Value = -304.5

$ws.Range("H27").Value = 1796.4333
$ws.Range("I27").Value = 599.5
$ws.Range("K27").Value = 599.5
$ws.Range("M27").Value = -492.5

$ws.Range("H28").Value = 5001
$ws.Range("J28").Value = 5001
$ws.Range("L28").Value = 5001
$ws.Range("N28").Value = -5465

$ws.Range("H37").Value = 5001
$ws.Range("J37").Value = 5001
$ws.Range("L37").Value = 5001
$ws.Range("N37").Value = -5215

$ws.Range("H136").Value = 10732.533
$ws.Range("I136").Value = 4333.1113
$ws.Range("J136").Value = 20331.666
$ws.Range("K136").Value = 12999.3339
$ws.Range("L136").Value = 60994.99800000001
$ws.Range("M136").Value = -10449.3339
$ws.Range("N136").Value = -66094.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 320
$ws.Range("J4").Value = 800
$ws.Range("L4").Value = 800
$ws.Range("N4").Value = -1026

$ws.Range("H14").Value = 3900
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 3900
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 3900
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -4236

